# Quarterly financial update for KMT.xlsx
# - Inserts two new columns (D, E) for the two most recent quarters (2018-09-30 and 2018-12-31)
# - Shifts the existing quarterly data two columns to the right (now F:M)
# - Populates the new D:E columns with the latest reported figures
# - Applies a handful of restatements to previously reported quarters that changed in this refresh

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KMT")

# Insert two blank columns before column D, across the populated data range (rows 5-102)
$ws.Range("D5:E102").EntireColumn.Insert()

# Copy number formatting from the (now shifted) old D:E columns -- now at F:G -- onto the new D:E columns
$ws.Range("F5:G102").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate the two new columns with the latest quarterly figures
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 587400
$ws.Cells.Item(8, 5).Value = 586700
$ws.Cells.Item(9, 4).Value = 388800
$ws.Cells.Item(9, 5).Value = 375600
$ws.Cells.Item(10, 4).Value = 198600
$ws.Cells.Item(10, 5).Value = 211100
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 1500
$ws.Cells.Item(14, 5).Value = 1100
$ws.Cells.Item(15, 4).Value = 3600
$ws.Cells.Item(15, 5).Value = 3600
$ws.Cells.Item(17, 4).Value = 508500
$ws.Cells.Item(17, 5).Value = 503500
$ws.Cells.Item(18, 4).Value = 78900
$ws.Cells.Item(18, 5).Value = 83200
$ws.Cells.Item(20, 4).Value = 4000
$ws.Cells.Item(20, 5).Value = 2700
$ws.Cells.Item(21, 4).Value = 110300
$ws.Cells.Item(21, 5).Value = 113500
$ws.Cells.Item(22, 4).Value = 8100
$ws.Cells.Item(22, 5).Value = 8100
$ws.Cells.Item(23, 4).Value = 74800
$ws.Cells.Item(23, 5).Value = 77800
$ws.Cells.Item(24, 4).Value = 22400
$ws.Cells.Item(24, 5).Value = 18400
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 52300
$ws.Cells.Item(26, 5).Value = 59400
$ws.Cells.Item(27, 4).Value = 50800
$ws.Cells.Item(27, 5).Value = 57700
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = 3900
$ws.Cells.Item(29, 5).Value = -1000
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = -4000
$ws.Cells.Item(32, 5).Value = -2700
$ws.Cells.Item(33, 4).Value = 54700
$ws.Cells.Item(33, 5).Value = 56700
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 54700
$ws.Cells.Item(35, 5).Value = 56700
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 96300
$ws.Cells.Item(41, 5).Value = 102100
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 4).Value = 380700
$ws.Cells.Item(43, 5).Value = 386700
$ws.Cells.Item(44, 4).Value = 578600
$ws.Cells.Item(44, 5).Value = 569300
$ws.Cells.Item(45, 4).Value = 63500
$ws.Cells.Item(45, 5).Value = 63500
$ws.Cells.Item(46, 4).Value = 1119000
$ws.Cells.Item(46, 5).Value = 1121500
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 4).Value = 855100
$ws.Cells.Item(48, 5).Value = 834400
$ws.Cells.Item(49, 4).Value = 468500
$ws.Cells.Item(49, 5).Value = 473300
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 88500
$ws.Cells.Item(52, 5).Value = 83100
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 2531100
$ws.Cells.Item(54, 5).Value = 2512200
$ws.Cells.Item(57, 4).Value = 198400
$ws.Cells.Item(57, 5).Value = 220900
$ws.Cells.Item(58, 4).Value = 3400
$ws.Cells.Item(58, 5).Value = 800
$ws.Cells.Item(59, 4).Value = 210300
$ws.Cells.Item(59, 5).Value = 217500
$ws.Cells.Item(60, 4).Value = 412100
$ws.Cells.Item(60, 5).Value = 439200
$ws.Cells.Item(61, 4).Value = 591700
$ws.Cells.Item(61, 5).Value = 591300
$ws.Cells.Item(62, 4).Value = 219100
$ws.Cells.Item(62, 5).Value = 217600
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 1261400
$ws.Cells.Item(66, 5).Value = 1284600
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 979300
$ws.Cells.Item(72, 5).Value = 941000
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 1269700
$ws.Cells.Item(76, 5).Value = 1227600
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 54700
$ws.Cells.Item(81, 5).Value = 56700
$ws.Cells.Item(83, 4).Value = 27400
$ws.Cells.Item(83, 5).Value = 27600
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 52300
$ws.Cells.Item(89, 5).Value = 9200
$ws.Cells.Item(91, 4).Value = -44800
$ws.Cells.Item(91, 5).Value = -43300
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -43100
$ws.Cells.Item(94, 5).Value = -42400
$ws.Cells.Item(96, 4).Value = -16400
$ws.Cells.Item(96, 5).Value = -16400
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = -13800
$ws.Cells.Item(100, 5).Value = -418900
$ws.Cells.Item(101, 4).Value = -1200
$ws.Cells.Item(101, 5).Value = -2000
$ws.Cells.Item(102, 4).Value = -5800
$ws.Cells.Item(102, 5).Value = -454100

# Restated figures for previously reported quarters that changed with this refresh
$ws.Cells.Item(9, 8).Value = 381800
$ws.Cells.Item(9, 9).Value = 716700
$ws.Cells.Item(10, 8).Value = 189500
$ws.Cells.Item(10, 9).Value = -174200
$ws.Cells.Item(14, 9).Value = 6800
$ws.Cells.Item(17, 8).Value = 507700
$ws.Cells.Item(17, 9).Value = 490300
$ws.Cells.Item(18, 8).Value = 63600
$ws.Cells.Item(18, 9).Value = 52200
$ws.Cells.Item(20, 8).Value = 3300
$ws.Cells.Item(20, 9).Value = 4200
$ws.Cells.Item(32, 8).Value = -3300
$ws.Cells.Item(32, 9).Value = -4200
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(47, 9).Value = 0
$ws.Cells.Item(47, 10).Value = 0
$ws.Cells.Item(91, 7).Value = -68800
$ws.Cells.Item(91, 8).Value = -37200
$ws.Cells.Item(91, 9).Value = -22300
$ws.Cells.Item(91, 10).Value = -23900
$ws.Cells.Item(94, 7).Value = -67400
$ws.Cells.Item(94, 8).Value = -36500
$ws.Cells.Item(94, 9).Value = -21900
